$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 0.07
$ws.Range("I3").Value = 0.04
$ws.Range("I4").Select()
